# Updated cryptos list on Fri Nov  8 03:45:10 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "75.760.21"
$ws.Range("E2").Value = "  +0.79%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.895.52"
$ws.Range("E3").Value = "  +1.02%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - Solana
$ws.Range("D5").Value = "197.55"
$ws.Range("E5").Value = "  +4.29%  "

# Row 6 - BNB
$ws.Range("D6").Value = "596.00"
$ws.Range("E6").Value = "  -1.76%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.01%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -3.02%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.80%  "

# Row 10 - LidoStakedEther
$ws.Range("D10").Value = "2.892.79"
$ws.Range("E10").Value = "  +1.05%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "0.417"
$ws.Range("E11").Value = "  +11.18%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -1.62%  "

# Row 13 - Toncoin
$ws.Range("E13").Value = "  -2.33%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.422.56"
$ws.Range("E14").Value = "  +1.02%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "75.604.13"
$ws.Range("E15").Value = "  +0.59%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -1.26%  "

# Row 17 - Avalanche
$ws.Range("D17").Value = "27.20"
$ws.Range("E17").Value = "  -2.12%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.889.36"
$ws.Range("E18").Value = "  +1.09%  "

# Row 19 - Uniswap
$ws.Range("E19").Value = "  -3.91%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "12.56"
$ws.Range("E20").Value = "  -0.65%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "375.55"
$ws.Range("E21").Value = "  -1.21%  "

# Row 22 - SuiNetwork
$ws.Range("E22").Value = "  -0.69%  "

# Row 23 - Polkadot
$ws.Range("E23").Value = "  -0.83%  "

# Row 24 - Dai
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.11%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "70.96"
$ws.Range("E25").Value = "  -0.65%  "

# Row 26 - WrappedeETH
$ws.Range("D26").Value = "3.043.27"
$ws.Range("E26").Value = "  +1.17%  "

# Row 27 - NEARProtocol
$ws.Range("E27").Value = "  -2.36%  "

# Row 28 - Aptos
$ws.Range("D28").Value = "9.51"
$ws.Range("E28").Value = "  -2.04%  "

# Row 29 - PEPE
$ws.Range("E29").Value = "  +2.12%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  -0.05%  "

# Row 31 - Fetch.AI
$ws.Range("E31").Value = "  -2.42%  "

# Row 32 - Bittensor
$ws.Range("D32").Value = "500.78"
$ws.Range("E32").Value = "  -6.27%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "7.69"
$ws.Range("E33").Value = "  -3.82%  "

# Row 34 - PancakeSwap
$ws.Range("E34").Value = "  -2.53%  "

# Row 35 - FirstDigitalUSD
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.02%  "

# Row 36 - Monero
$ws.Range("D36").Value = "163.32"
$ws.Range("E36").Value = "  +0.47%  "

# Row 37 - EthereumClassic
$ws.Range("E37").Value = "  -2.26%  "

# Row 38 - WhiteBITCoin
$ws.Range("D38").Value = "19.70"
$ws.Range("E38").Value = "  +2.04%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -7.31%  "

# Row 40 - USDe
$ws.Range("E40").Value = "  -0.12%  "

# Row 41 - Aave
$ws.Range("D41").Value = "178.90"
$ws.Range("E41").Value = "  -3.07%  "

# Row 42 - PolygonEcosystemToken
$ws.Range("E42").Value = "  -1.66%  "

# Row 43 - RenderToken
$ws.Range("E43").Value = "  -3.78%  "

# Row 44 - Stacks
$ws.Range("D44").Value = "1.66"
$ws.Range("E44").Value = "  -3.10%  "

# Row 45 - Cronos
$ws.Range("D45").Value = "0.0902"
$ws.Range("E45").Value = "  +4.71%  "

# Row 46 - now ImmutableX (was OKB)
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").Value = "1.20"
$ws.Range("E46").Value = "  -5.73%  "

# Row 47 - now OKB (was ImmutableX)
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "40.03"
$ws.Range("E47").Value = "  -0.05%  "

# Row 48 - dogwifhat
$ws.Range("D48").Value = "2.33"
$ws.Range("E48").Value = "  -2.62%  "

# Row 49 - ARBITRUM
$ws.Range("D49").Value = "0.574"
$ws.Range("E49").Value = "  -1.34%  "

# Row 50 - Filecoin
$ws.Range("E50").Value = "  -2.30%  "

# Row 51 - Mantle
$ws.Range("E51").Value = "  +4.83%  "
